$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new 2022 column (S), mirroring column R's formatting ---
# Pull the formats (number format / font / border / alignment) from R3:R4
# onto S3:S4 the same way a user would via copy + paste-special-formats,
# then fill in the new year label and indicator value.
$ws.Range("R3:R4").Copy()
$ws.Range("S3:S4").PasteSpecial(-4122)

$ws.Range("S3").Value = 2022
$ws.Range("S4").Value = 0.071025550219041236

# --- Re-style the first three columns to one uniform width ---
$ws.Range("A1:C1").ColumnWidth = 32.666666666666664

# --- Move the saved selection to match the new view state ---
$ws.Range("F14").Select()
